# Common: Initial boring stuff
# Adds a new "Translations - Market" sheet (mirroring the other
# "Translations - *" sheets), registers it in the "tabs" index sheet,
# and leaves the view state (selections / scroll position / active
# sheet) the way the author left it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Register the new tab in the "tabs" sheet (row 9).
# ---------------------------------------------------------------
$tabs = $wb.Worksheets.Item("tabs")
$tabs.Range("A9").Value = "Translations - Market"
$tabs.Range("B9").Value = "translation"

# ---------------------------------------------------------------
# 2. Create the new worksheet after "Translations - Shared".
# ---------------------------------------------------------------
$shared = $wb.Worksheets.Item("Translations - Shared")
$market = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $shared)
$market.Name = "Translations - Market"

# Column widths (characters) matching the other translation sheets.
$market.Columns.Item(1).ColumnWidth = 25.83
$market.Columns.Item(2).ColumnWidth = 37.67
$market.Columns.Item(3).ColumnWidth = 35

# Header row.
$market.Range("A1").Value = "Language"
$market.Range("B1").Value = "Label"
$market.Range("C1").Value = "Text"
$market.Range("A1:C1").Style = "Nadpis 2"

# Data rows.
$rows = @(
    @("cs", "market.index.title", "Tržiště"),
    @("cs", "market.home.subtitle", "Tato část aplikace slouží k získávání předmětů, které již fyzicky vlastníte, ale přejete si je zaevidovat do systému."),
    @("cs", "market.home.menu", "Domů"),
    @("cs", "market.lab.menu", "Laboratoř"),
    @("cs", "market.home.title", "Tržiště"),
    @("cs", "market.atomizer.menu", "Atomizéry"),
    @("cs", "market.mod.menu", "Mody"),
    @("cs", "market.cotton.menu", "Vaty"),
    @("cs", "market.cell.menu", "Články"),
    @("cs", "market.atomizer.index.title", "Atomizéry"),
    @("cs", "market.mod.index.title", "Mody"),
    @("cs", "market.cotton.index.title", "Vaty"),
    @("cs", "market.cell.index.title", "Články")
)

$r = 2
foreach ($row in $rows) {
    $market.Range("A$r").Value = $row[0]
    $market.Range("B$r").Value = $row[1]
    $market.Range("C$r").Value = $row[2]
    $market.Range("A$r`:C$r").Style = "import"
    $r++
}

# Row 3 holds the long subtitle translation and wraps to three lines.
$market.Rows.Item(3).RowHeight = 39

$market.PageSetup.PaperSize = 9
$market.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 3. Restore the view state on the sheets whose scroll / selection
#    moved around while the author was working.
# ---------------------------------------------------------------

# "Translations - Public" keeps its scroll position but is no longer
# the active tab.
$public = $wb.Worksheets.Item("Translations - Public")
$public.Range("C20").Select()

# "Translations - Shared" scrolled back to the top and selection
# moved to C9.
$shared.Activate()
$shared.Range("C9").Select()

# "tabs" selection moved to the freshly added row.
$tabs.Activate()
$tabs.Range("A9").Select()

# Finally, the new sheet is the active tab, with B7 selected.
$market.Activate()
$market.Range("B7").Select()
